$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.120.38'
$ws.Range("E2").Value = '  -3.55%  '

$ws.Range("D3").Value = '1.605.61'
$ws.Range("E3").Value = '  -2.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.37'
$ws.Range("E6").Value = '  -2.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3767'
$ws.Range("E7").Value = '  -3.59%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3657'
$ws.Range("E8").Value = '  -4.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.63'
$ws.Range("E9").Value = '  -5.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9991'
$ws.Range("E10").Value = '  -0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.277'
$ws.Range("E11").Value = '  -5.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08085'
$ws.Range("E12").Value = '  -4.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.99'
$ws.Range("E13").Value = '  -4.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.632'
$ws.Range("E14").Value = '  -6.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.677'
$ws.Range("E15").Value = '  -2.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001266'
$ws.Range("E16").Value = '  -3.75%  '

$ws.Range("D17").Value = '1.597.76'
$ws.Range("E17").Value = '  -3.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.40'
$ws.Range("E18").Value = '  -3.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06787'
$ws.Range("E19").Value = '  -2.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.39'
$ws.Range("E20").Value = '  -6.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.580'
$ws.Range("E21").Value = '  -4.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.01'
$ws.Range("E23").Value = '  -4.45%  '

$ws.Range("D24").Value = '23.127.45'
$ws.Range("E24").Value = '  -3.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.344'
$ws.Range("E25").Value = '  -5.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.928'
$ws.Range("E26").Value = '  -2.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.45'
$ws.Range("E28").Value = '  -1.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.277'
$ws.Range("E29").Value = '  -3.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.06'
$ws.Range("E30").Value = '  -5.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.402'
$ws.Range("E31").Value = '  -3.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.950'
$ws.Range("E32").Value = '  -10.29%  '

$ws.Range("D33").Value = '1.769.92'
$ws.Range("E33").Value = '  -3.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9939'
$ws.Range("E34").Value = '  -4.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07720'
$ws.Range("E35").Value = '  -4.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02787'
$ws.Range("E36").Value = '  -6.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.292'
$ws.Range("E37").Value = '  -6.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2543'
$ws.Range("E38").Value = '  -5.20%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.11'
$ws.Range("E39").Value = '  -6.72%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08846'
$ws.Range("E40").Value = '  -3.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.396'
$ws.Range("E41").Value = '  -1.91%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7160'
$ws.Range("E42").Value = '  -5.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.78'
$ws.Range("E43").Value = '  -5.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.80'
$ws.Range("E44").Value = '  -2.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6612'
$ws.Range("E45").Value = '  -4.84%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.308'
$ws.Range("E46").Value = '  -6.03%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.966'
$ws.Range("E48").Value = '  -2.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08005'
$ws.Range("E49").Value = '  -3.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.26'
$ws.Range("E50").Value = '  -2.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.171'
$ws.Range("E51").Value = '  -4.25%  '
